$d = $word.ActiveDocument

# The floating "Classification: Controlled" text box lives in the default
# (primary) footer of the document's only section. Locate it via the
# Section's Footers collection and remove the shape entirely, leaving the
# empty Footer-styled paragraph behind (matching the target OOXML).
$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)

for ($i = $footer.Shapes.Count; $i -ge 1; $i--) {
    $shape = $footer.Shapes.Item($i)
    $shape.Delete()
}
